$d = $word.ActiveDocument

# --- 1. The cursor used to sit (marked by the hidden "_GoBack" bookmark)
#        right in the middle of "...subroutine if ther|e was a 0...", which
#        had split that sentence across two runs. Re-typing over that exact
#        span collapses it back into a single run (and removes the
#        now-stale bookmark sitting inside it) without altering the text.

$rng = $d.Content
$rng.Find.Execute("subroutine if there was a 0", $true, $false, $false, $false, $false, $true, 1, $false, "subroutine if there was a 0", 2) | Out-Null

# --- 2. Insert a brand-new paragraph right after the existing final
#        paragraph (the "timer drivers" paragraph), before the section
#        break, and fill it with the new "Finally, we implemented 2 simple
#        timers..." text describing the stopwatch/pushbutton polling.

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter() | Out-Null

$d2 = $word.ActiveDocument
$newLast = $d2.Paragraphs.Last
$newLast.Range.Text = "Finally, we implemented 2 simple timers with this driver, with a 1 second and a 5ms timeout respectively. The actual stopwatch updated at every 10 milliseconds and rewrote the displays then too. We could maybe have improved the logic by only rewriting if there was a change to the value, but this was the most straightforward and error free way to implement it. This functionality was governed by the polling of the pushbuttons, where we either started, stopped or reset the timer depending on the pushbuttons pressed. This could have been improved by using the last pushbutton to implement increased functionality, such as loading in a time to start from using the switches from the previous section, perhaps even counting backwards from there. However, the conversion of the switch values into time values would have been daunting to even consider, and we did not have time to think about it."

# --- 3. Move the "_GoBack" bookmark (Word's last-edit-position marker) so
#        it again sits at the very end of the document, i.e. at the end of
#        the paragraph we just typed into (re-adding a bookmark with the
#        same name relocates it instead of creating a duplicate).

$d3 = $word.ActiveDocument
$finalRange = $d3.Paragraphs.Last.Range
$endSpot = $d3.Range($finalRange.End - 2, $finalRange.End - 1)
$d3.Bookmarks.Add("_GoBack", $endSpot) | Out-Null

Write-Output "done"
